# Refresh the cryptos price/volume table (GitHub Actions daily scrape update).
# For Price (column D) cells whose new value is purely numeric-looking
# (e.g. "607.46"), force text storage via NumberFormat="@" before assigning,
# then reset the display style back to "Normal" so no stray style index is
# left on the cell - matching the original inline-string/no-style cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.749.56'
$ws.Range('D3').Value = '3.808.86'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.451'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000253'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.06'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '4.445.07'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').Value = '3.814.05'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '67.777.50'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '462.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.30%  '
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  -3.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('E26').Value = '  -1.36%  '
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').Value = '3.958.00'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('E31').Value = '  +1.33%  '
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('D36').Value = '3.749.12'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  +2.37%  '
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.27%  '
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '148.80'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.37'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.42%  '
$ws.Range('E51').Value = '  +0.42%  '
